# Commit: تعديل تلقائي في شيت Card12 by admin at 2025-12-06 18:30:23
# Update the "card" column (A) on sheet Card12 from "2" to "12" for the
# service-band rows that still carried the old value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Rows 3,4,5,6,7,9,10,11,12,13 currently hold "2" in column A and need to
# become "12" (rows 2 and 8 already read "12" and are left untouched).
$rows = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 13)
foreach ($r in $rows) {
    # Leading apostrophe keeps the value text (matching the existing
    # "card" column entries) instead of letting Excel coerce it to a number.
    $ws.Cells.Item($r, 1).Value = "'12"
}
